$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.06437833333333333
$ws.Range("H2").Value = 0.193135
$ws.Range("I2").Value = 0.109187438766332
$ws.Range("J2").Value = 0.109187438766332
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 9.379971107065556
$ws.Range("R2").Value = 84.41973996359002
$ws.Range("S2").Value = 0.03129244025978011
$ws.Range("T2").Value = 0.03129244025978011
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.06437833333333333
$ws.Range("H3").Value = 0.193135
$ws.Range("I3").Value = 0.109187438766332
$ws.Range("J3").Value = 0.109187438766332
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 10.86704352484222
$ws.Range("R3").Value = 97.80339172358002
$ws.Range("S3").Value = 0.03625344965566096
$ws.Range("T3").Value = 0.03625344965566096
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.06437833333333333
$ws.Range("H4").Value = 0.193135
$ws.Range("I4").Value = 0.109187438766332
$ws.Range("J4").Value = 0.109187438766332
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 8.24854829384889
$ws.Range("R4").Value = 74.23693464464
$ws.Range("S4").Value = 0.02751791042519824
$ws.Range("T4").Value = 0.02751791042519824
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.06437833333333333
$ws.Range("H5").Value = 0.193135
$ws.Range("I5").Value = 0.109187438766332
$ws.Range("J5").Value = 0.109187438766332
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 4.233588664221666
$ws.Range("R5").Value = 38.102297977995
$ws.Range("S5").Value = 0.01412363842569274
$ws.Range("T5").Value = 0.01412363842569274
$ws.Range("G6").Value = 0.4788196666666666
$ws.Range("I6").Value = 0.8120914339857952
$ws.Range("J6").Value = 0.8120914339857951
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 69.76438199437845
$ws.Range("R6").Value = 627.8794379494061
$ws.Range("S6").Value = 0.2327403497197477
$ws.Range("T6").Value = 0.2327403497197477
$ws.Range("G7").Value = 0.4788196666666666
$ws.Range("I7").Value = 0.8120914339857952
$ws.Range("J7").Value = 0.8120914339857951
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("Q7").Value = 80.82461736428579
$ws.Range("R7").Value = 727.4215562785721
$ws.Range("S7").Value = 0.2696383050142185
$ws.Range("T7").Value = 0.2696383050142184
$ws.Range("G8").Value = 0.4788196666666666
$ws.Range("I8").Value = 0.8120914339857952
$ws.Range("J8").Value = 0.8120914339857951
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 61.34932266877512
$ws.Range("R8").Value = 552.1439040189759
$ws.Range("S8").Value = 0.2046669432856284
$ws.Range("T8").Value = 0.2046669432856283
$ws.Range("G9").Value = 0.4788196666666666
$ws.Range("I9").Value = 0.8120914339857952
$ws.Range("J9").Value = 0.8120914339857951
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 31.48769792642033
$ws.Range("R9").Value = 283.389281337783
$ws.Range("S9").Value = 0.1050458359662007
$ws.Range("T9").Value = 0.1050458359662007
$ws.Range("G10").Value = 0.042481
$ws.Range("H10").Value = 0.127443
$ws.Range("I10").Value = 0.07204895414449818
$ws.Range("J10").Value = 0.07204895414449818
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 6.189513334184667
$ws.Range("R10").Value = 55.70562000766201
$ws.Range("S10").Value = 0.02064878175383621
$ws.Range("T10").Value = 0.02064878175383621
$ws.Range("G11").Value = 0.042481
$ws.Range("H11").Value = 0.127443
$ws.Range("I11").Value = 0.07204895414449818
$ws.Range("J11").Value = 0.07204895414449818
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 7.170780168982668
$ws.Range("R11").Value = 64.53702152084401
$ws.Range("S11").Value = 0.02392237753108655
$ws.Range("T11").Value = 0.02392237753108655
$ws.Range("G12").Value = 0.042481
$ws.Range("H12").Value = 0.127443
$ws.Range("I12").Value = 0.07204895414449818
$ws.Range("J12").Value = 0.07204895414449818
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 5.442927176394667
$ws.Range("R12").Value = 48.986344587552
$ws.Range("S12").Value = 0.0181581021478165
$ws.Range("T12").Value = 0.0181581021478165
$ws.Range("G13").Value = 0.042481
$ws.Range("H13").Value = 0.127443
$ws.Range("I13").Value = 0.07204895414449818
$ws.Range("J13").Value = 0.07204895414449818
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 2.793596396999
$ws.Range("R13").Value = 25.142367572991
$ws.Range("S13").Value = 0.009319692711758925
$ws.Range("T13").Value = 0.009319692711758925
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.003934
$ws.Range("H14").Value = 0.011802
$ws.Range("I14").Value = 0.006672173103374587
$ws.Range("J14").Value = 0.006672173103374586
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 0.5731867295186667
$ws.Range("R14").Value = 5.158680565668001
$ws.Range("S14").Value = 0.001912203277220208
$ws.Range("T14").Value = 0.001912203277220208
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.003934
$ws.Range("H15").Value = 0.011802
$ws.Range("I15").Value = 0.006672173103374587
$ws.Range("J15").Value = 0.006672173103374586
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 0.6640580302906668
$ws.Range("R15").Value = 5.976522272616001
$ws.Range("S15").Value = 0.002215358235618146
$ws.Range("T15").Value = 0.002215358235618146
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.003934
$ws.Range("H16").Value = 0.011802
$ws.Range("I16").Value = 0.006672173103374587
$ws.Range("J16").Value = 0.006672173103374586
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 0.5040482924586667
$ws.Range("R16").Value = 4.536434632128
$ws.Range("S16").Value = 0.001681551136967352
$ws.Range("T16").Value = 0.001681551136967352
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.003934
$ws.Range("H17").Value = 0.011802
$ws.Range("I17").Value = 0.006672173103374587
$ws.Range("J17").Value = 0.006672173103374586
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 0.258704084786
$ws.Range("R17").Value = 2.328336763074
$ws.Range("S17").Value = 0.0008630604535688805
$ws.Range("T17").Value = 0.0008630604535688804
